# Rename the worksheet "Property1" to "DataNode" to unify the
# conception of DataNode / DataTable / Entity across the data configs.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Property1")
$ws.Name = "DataNode"

# Preserve the author's final cell selection recorded in the commit
# (bottom pane, cell O40) so the saved view state matches too.
[void]$ws.Range("O40").Select()
